$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update res_bus vm_pu results after changing slack bus voltage setpoint from 1.05 to 1.02 (380 kV case)
# Row 2
$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.048787380830204
$ws.Cells.Item(2, 4).Value = 1.05578431934836
$ws.Cells.Item(2, 5).Value = 0.992614727750844
$ws.Cells.Item(2, 6).Value = 1.064452307921284
$ws.Cells.Item(2, 9).Value = 1.046655076873322
$ws.Cells.Item(2, 10).Value = 1.053829278645595
$ws.Cells.Item(2, 11).Value = 1.058523825643732
$ws.Cells.Item(2, 12).Value = 0.9955398523335997
$ws.Cells.Item(2, 13).Value = 1.067168206705943
$ws.Cells.Item(2, 14).Value = 1.021752181434427

# Row 3
$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.049751318248072
$ws.Cells.Item(3, 4).Value = 1.056557194739944
$ws.Cells.Item(3, 5).Value = 0.9936372048519299
$ws.Cells.Item(3, 6).Value = 1.065356026627811
$ws.Cells.Item(3, 9).Value = 1.046919442814507
$ws.Cells.Item(3, 10).Value = 1.054442163160081
$ws.Cells.Item(3, 11).Value = 1.059110314433211
$ws.Cells.Item(3, 12).Value = 0.9963617723202687
$ws.Cells.Item(3, 13).Value = 1.067886911492349
$ws.Cells.Item(3, 14).Value = 1.021958928704951

# Row 4
$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.050375327028714
$ws.Cells.Item(4, 4).Value = 1.05705755323334
$ws.Cells.Item(4, 5).Value = 0.9942998659930998
$ws.Cells.Item(4, 6).Value = 1.065941379336871
$ws.Cells.Item(4, 9).Value = 1.047089509412093
$ws.Cells.Item(4, 10).Value = 1.054838384406281
$ws.Cells.Item(4, 11).Value = 1.059489406525065
$ws.Cells.Item(4, 12).Value = 0.9968940712668347
$ws.Cells.Item(4, 13).Value = 1.068351911154328
$ws.Cells.Item(4, 14).Value = 1.022092514844582

# Row 5
$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.05063772582993
$ws.Cells.Item(5, 4).Value = 1.057267964168368
$ws.Cells.Item(5, 5).Value = 0.994578699834602
$ws.Cells.Item(5, 6).Value = 1.066187600759304
$ws.Cells.Item(5, 9).Value = 1.047160766537858
$ws.Cells.Item(5, 10).Value = 1.055004869796042
$ws.Cells.Item(5, 11).Value = 1.059648678943662
$ws.Cells.Item(5, 12).Value = 0.9971179600053012
$ws.Cells.Item(5, 13).Value = 1.068547383946813
$ws.Cells.Item(5, 14).Value = 1.022148627876382

# Row 6
$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.050681787562794
$ws.Cells.Item(6, 4).Value = 1.057303296587896
$ws.Cells.Item(6, 5).Value = 0.994625531979634
$ws.Cells.Item(6, 6).Value = 1.066228950543793
$ws.Cells.Item(6, 9).Value = 1.047172716901791
$ws.Cells.Item(6, 10).Value = 1.055032818358219
$ws.Cells.Item(6, 11).Value = 1.059675415718184
$ws.Cells.Item(6, 12).Value = 0.9971555583673455
$ws.Cells.Item(6, 13).Value = 1.068580203893558
$ws.Cells.Item(6, 14).Value = 1.022158046760009

# Row 7
$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.050378832957869
$ws.Cells.Item(7, 4).Value = 1.057060364519373
$ws.Cells.Item(7, 5).Value = 0.994303590798249
$ws.Cells.Item(7, 6).Value = 1.065944668814345
$ws.Cells.Item(7, 9).Value = 1.047090462492474
$ws.Cells.Item(7, 10).Value = 1.054840609331201
$ws.Cells.Item(7, 11).Value = 1.059491535116646
$ws.Cells.Item(7, 12).Value = 0.9968970624462089
$ws.Cells.Item(7, 13).Value = 1.068354523123938
$ws.Cells.Item(7, 14).Value = 1.022093264813112

# Row 8
$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.049113090219373
$ws.Cells.Item(8, 4).Value = 1.056045462407878
$ws.Cells.Item(8, 5).Value = 0.9929600610674297
$ws.Cells.Item(8, 6).Value = 1.064757601768098
$ws.Cells.Item(8, 9).Value = 1.046744626497049
$ws.Cells.Item(8, 10).Value = 1.054036479076625
$ws.Cells.Item(8, 11).Value = 1.058722115699898
$ws.Cells.Item(8, 12).Value = 0.9958175282591056
$ws.Cells.Item(8, 13).Value = 1.067411106185572
$ws.Cells.Item(8, 14).Value = 1.021822092437374

# Row 9
$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.046884843321722
$ws.Cells.Item(9, 4).Value = 1.054259090605662
$ws.Cells.Item(9, 5).Value = 0.9906006454969559
$ws.Cells.Item(9, 6).Value = 1.062670381577712
$ws.Cells.Item(9, 9).Value = 1.046127612827587
$ws.Cells.Item(9, 10).Value = 1.052616806999259
$ws.Cells.Item(9, 11).Value = 1.05736323316295
$ws.Cells.Item(9, 12).Value = 0.9939188001724441
$ws.Cells.Item(9, 13).Value = 1.065748343177962
$ws.Cells.Item(9, 14).Value = 1.021342787417987

# Row 10
$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.04540083668444
$ws.Cells.Item(10, 4).Value = 1.053069599453732
$ws.Cells.Item(10, 5).Value = 0.989033133672735
$ws.Cells.Item(10, 6).Value = 1.061282027257948
$ws.Cells.Item(10, 9).Value = 1.045711183339656
$ws.Cells.Item(10, 10).Value = 1.05166859504695
$ws.Cells.Item(10, 11).Value = 1.056455298806346
$ws.Cells.Item(10, 12).Value = 0.9926553831429383
$ws.Cells.Item(10, 13).Value = 1.064639663153154
$ws.Cells.Item(10, 14).Value = 1.021022285601294

# Row 11
$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.044758607723088
$ws.Cells.Item(11, 4).Value = 1.052554888364296
$ws.Cells.Item(11, 5).Value = 0.988355674866747
$ws.Cells.Item(11, 6).Value = 1.060681610512617
$ws.Cells.Item(11, 9).Value = 1.045529663718284
$ws.Cells.Item(11, 10).Value = 1.051257599770938
$ws.Cells.Item(11, 11).Value = 1.056061685969938
$ws.Cells.Item(11, 12).Value = 0.9921088820399291
$ws.Cells.Item(11, 13).Value = 1.064159563814273
$ws.Cells.Item(11, 14).Value = 1.020883279880055

# Row 12
$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.044520109351828
$ws.Cells.Item(12, 4).Value = 1.052363754894417
$ws.Cells.Item(12, 5).Value = 0.9881042295826724
$ws.Cells.Item(12, 6).Value = 1.06045870279038
$ws.Cells.Item(12, 9).Value = 1.045462058933614
$ws.Cells.Item(12, 10).Value = 1.051104876444603
$ws.Cells.Item(12, 11).Value = 1.055915410606381
$ws.Cells.Item(12, 12).Value = 0.9919059725120875
$ws.Cells.Item(12, 13).Value = 1.063981229421466
$ws.Cells.Item(12, 14).Value = 1.020831613271626

# Row 13
$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.044571265605734
$ws.Cells.Item(13, 4).Value = 1.052404751249896
$ws.Cells.Item(13, 5).Value = 0.9881581567098651
$ws.Cells.Item(13, 6).Value = 1.060506512085968
$ws.Cells.Item(13, 9).Value = 1.045476568542485
$ws.Cells.Item(13, 10).Value = 1.051137638882845
$ws.Cells.Item(13, 11).Value = 1.055946790325022
$ws.Cells.Item(13, 12).Value = 0.9919494934313052
$ws.Cells.Item(13, 13).Value = 1.064019482925805
$ws.Cells.Item(13, 14).Value = 1.020842697454321

# Row 14
$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.044738892266906
$ws.Cells.Item(14, 4).Value = 1.052539088121526
$ws.Cells.Item(14, 5).Value = 0.9883348863814464
$ws.Cells.Item(14, 6).Value = 1.060663182562893
$ws.Cells.Item(14, 9).Value = 1.0455240791617
$ws.Cells.Item(14, 10).Value = 1.05124497686122
$ws.Cells.Item(14, 11).Value = 1.056049596226884
$ws.Cells.Item(14, 12).Value = 0.9920921077337197
$ws.Cells.Item(14, 13).Value = 1.06414482271648
$ws.Cells.Item(14, 14).Value = 1.020879009787659

# Row 15
$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.044842179767352
$ws.Cells.Item(15, 4).Value = 1.052621864577089
$ws.Cells.Item(15, 5).Value = 0.9884438009545853
$ws.Cells.Item(15, 6).Value = 1.060759727526856
$ws.Cells.Item(15, 9).Value = 1.045553328140982
$ws.Cells.Item(15, 10).Value = 1.051311103219775
$ws.Cells.Item(15, 11).Value = 1.05611292907945
$ws.Cells.Item(15, 12).Value = 0.9921799884222134
$ws.Cells.Item(15, 13).Value = 1.064222048179657
$ws.Cells.Item(15, 14).Value = 1.020901378559064

# Row 16
$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.045443466852555
$ws.Cells.Item(16, 4).Value = 1.053103766526264
$ws.Cells.Item(16, 5).Value = 0.9890781214508737
$ws.Cells.Item(16, 6).Value = 1.061321890835058
$ws.Cells.Item(16, 9).Value = 1.045723204890031
$ws.Cells.Item(16, 10).Value = 1.051695862803316
$ws.Cells.Item(16, 11).Value = 1.056481411730766
$ws.Cells.Item(16, 12).Value = 0.9926916645766087
$ws.Cells.Item(16, 13).Value = 1.064671525136444
$ws.Cells.Item(16, 14).Value = 1.021031506211174

# Row 17
$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.04582073402379
$ws.Cells.Item(17, 4).Value = 1.053406144412184
$ws.Cells.Item(17, 5).Value = 0.989476357848556
$ws.Cells.Item(17, 6).Value = 1.061674722580481
$ws.Cells.Item(17, 9).Value = 1.045829442321945
$ws.Cells.Item(17, 10).Value = 1.051937102184269
$ws.Cells.Item(17, 11).Value = 1.056712425535579
$ws.Cells.Item(17, 12).Value = 0.9930127773699352
$ws.Cells.Item(17, 13).Value = 1.064953461818686
$ws.Cells.Item(17, 14).Value = 1.021113071474426

# Row 18
$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.046040821860577
$ws.Cells.Item(18, 4).Value = 1.053582549668868
$ws.Cells.Item(18, 5).Value = 0.9897087662937556
$ws.Cells.Item(18, 6).Value = 1.061880595640313
$ws.Cells.Item(18, 9).Value = 1.045891292728223
$ws.Cells.Item(18, 10).Value = 1.052077773114446
$ws.Cells.Item(18, 11).Value = 1.056847126457047
$ws.Cells.Item(18, 12).Value = 0.9932001317071769
$ws.Cells.Item(18, 13).Value = 1.065117907414356
$ws.Cells.Item(18, 14).Value = 1.021160625255322

# Row 19
$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.04611587192081
$ws.Cells.Item(19, 4).Value = 1.053642704926249
$ws.Cells.Item(19, 5).Value = 0.9897880325774034
$ws.Cells.Item(19, 6).Value = 1.06195080525754
$ws.Cells.Item(19, 9).Value = 1.045912362415094
$ws.Cells.Item(19, 10).Value = 1.052125731473978
$ws.Cells.Item(19, 11).Value = 1.056893048245413
$ws.Cells.Item(19, 12).Value = 0.9932640239640975
$ws.Cells.Item(19, 13).Value = 1.065173978524872
$ws.Cells.Item(19, 14).Value = 1.021176836148988

# Row 20
$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.045780253276603
$ws.Cells.Item(20, 4).Value = 1.053373698669012
$ws.Cells.Item(20, 5).Value = 0.9894336180360679
$ws.Cells.Item(20, 6).Value = 1.061636859589674
$ws.Cells.Item(20, 9).Value = 1.045818056055732
$ws.Cells.Item(20, 10).Value = 1.05191122360568
$ws.Cells.Item(20, 11).Value = 1.056687644643606
$ws.Cells.Item(20, 12).Value = 0.9929783193494215
$ws.Cells.Item(20, 13).Value = 1.064923213015659
$ws.Cells.Item(20, 14).Value = 1.021104322552163

# Row 21
$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.044689528895851
$ws.Cells.Item(21, 4).Value = 1.052499527785165
$ws.Cells.Item(21, 5).Value = 0.9882828385668249
$ws.Cells.Item(21, 6).Value = 1.060617043862044
$ws.Cells.Item(21, 9).Value = 1.045510093442807
$ws.Cells.Item(21, 10).Value = 1.051213370190679
$ws.Cells.Item(21, 11).Value = 1.056019324378131
$ws.Cells.Item(21, 12).Value = 0.9920501090198102
$ws.Cells.Item(21, 13).Value = 1.064107913378514
$ws.Cells.Item(21, 14).Value = 1.02086831763821

# Row 22
$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.04400406027289
$ws.Cells.Item(22, 4).Value = 1.051950209886767
$ws.Cells.Item(22, 5).Value = 0.9875604150241495
$ws.Cells.Item(22, 6).Value = 1.0599765041578
$ws.Cells.Item(22, 9).Value = 1.045315422080884
$ws.Cells.Item(22, 10).Value = 1.050774247267753
$ws.Cells.Item(22, 11).Value = 1.055598720110158
$ws.Cells.Item(22, 12).Value = 0.9914670000341481
$ws.Cells.Item(22, 13).Value = 1.063595278746721
$ws.Cells.Item(22, 14).Value = 1.020719737189091

# Row 23
$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.044367410274424
$ws.Cells.Item(23, 4).Value = 1.052241384133216
$ws.Cells.Item(23, 5).Value = 0.9879432794643023
$ws.Cells.Item(23, 6).Value = 1.060316003536825
$ws.Cells.Item(23, 9).Value = 1.045418719784918
$ws.Cells.Item(23, 10).Value = 1.051007068031507
$ws.Cells.Item(23, 11).Value = 1.055821728502023
$ws.Cells.Item(23, 12).Value = 0.991776070289318
$ws.Cells.Item(23, 13).Value = 1.063867037973391
$ws.Cells.Item(23, 14).Value = 1.020798520857783

# Row 24
$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.045798544677864
$ws.Cells.Item(24, 4).Value = 1.053388359400785
$ws.Cells.Item(24, 5).Value = 0.9894529299347244
$ws.Cells.Item(24, 6).Value = 1.061653968022211
$ws.Cells.Item(24, 9).Value = 1.045823201378038
$ws.Cells.Item(24, 10).Value = 1.051922917144796
$ws.Cells.Item(24, 11).Value = 1.056698842202715
$ws.Cells.Item(24, 12).Value = 0.9929938892766442
$ws.Cells.Item(24, 13).Value = 1.064936881157969
$ws.Cells.Item(24, 14).Value = 1.021108275881099

# Row 25
$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.047460638895348
$ws.Cells.Item(25, 4).Value = 1.054720664657867
$ws.Cells.Item(25, 5).Value = 0.9912096547607049
$ws.Cells.Item(25, 6).Value = 1.063209432278699
$ws.Cells.Item(25, 9).Value = 1.046288024269205
$ws.Cells.Item(25, 10).Value = 1.052984140528175
$ws.Cells.Item(25, 11).Value = 1.057714895249285
$ws.Cells.Item(25, 12).Value = 0.9944092447426414
$ws.Cells.Item(25, 13).Value = 1.066178241758321
$ws.Cells.Item(25, 14).Value = 1.021466870801491
